$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns containing percentage-like text values ("NN%") must have their
# NumberFormat forced to Text ("@") before assignment, otherwise Excel
# auto-converts the literal "NN%" into a numeric percentage value.

$ws.Range("E2").Value = "2026-02-24 20:18:41"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "41%"
$ws.Range("O2").Value = "5.9 °C"
$ws.Range("E3").Value = "2026-02-24 20:18:43"
$ws.Range("E4").Value = "2026-02-24 20:18:46"
$ws.Range("O4").Value = "13.1 °C"
$ws.Range("E5").Value = "2026-02-24 20:18:49"
$ws.Range("E6").Value = "2026-02-24 20:18:52"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "68%"
$ws.Range("J6").Value = "1019.6 hPa"
$ws.Range("E7").Value = "2026-02-24 20:18:54"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "70%"
$ws.Range("J7").Value = "1020.2 hPa"
$ws.Range("O7").Value = "14.2 °C"
$ws.Range("E8").Value = "2026-02-24 20:18:56"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "46%"
$ws.Range("E9").Value = "2026-02-24 20:18:59"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "81%"
$ws.Range("E10").Value = "2026-02-24 20:19:01"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "76%"
$ws.Range("O10").Value = "11.4 °C"
$ws.Range("E11").Value = "2026-02-24 20:19:04"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "69%"
$ws.Range("E12").Value = "2026-02-24 20:19:07"
$ws.Range("E13").Value = "2026-02-24 20:19:09"
$ws.Range("O13").Value = "6.8 °C"
$ws.Range("E14").Value = "2026-02-24 20:19:12"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "84%"
$ws.Range("O14").Value = "11.7 °C"
$ws.Range("E15").Value = "2026-02-24 20:19:15"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "76%"
$ws.Range("O15").Value = "12.1 °C"
$ws.Range("E16").Value = "2026-02-24 20:19:17"
$ws.Range("G16").Value = "68 cm"
$ws.Range("L16").Value = "24.8 km/h - 184º 19:49 TU"
$ws.Range("E17").Value = "2026-02-24 20:19:19"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "33%"
$ws.Range("E18").Value = "2026-02-24 20:19:22"
$ws.Range("J18").Value = "1020.1 hPa"
$ws.Range("E19").Value = "2026-02-24 20:19:25"
$ws.Range("O19").Value = "12.7 °C"
$ws.Range("E20").Value = "2026-02-24 20:19:28"
$ws.Range("E21").Value = "2026-02-24 20:19:31"
$ws.Range("E22").Value = "2026-02-24 20:19:33"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "24%"
$ws.Range("N22").Value = "1.2 °C 19:54 TU"
$ws.Range("E23").Value = "2026-02-24 20:19:36"
$ws.Range("E24").Value = "2026-02-24 20:19:39"
$ws.Range("L24").Value = "10.4 km/h - 106º 19:48 TU"
$ws.Range("E25").Value = "2026-02-24 20:19:42"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = "33%"
$ws.Range("O25").Value = "6.8 °C"
$ws.Range("E26").Value = "2026-02-24 20:19:44"
$ws.Range("E27").Value = "2026-02-24 20:19:47"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "32%"
$ws.Range("O27").Value = "6.4 °C"
$ws.Range("E28").Value = "2026-02-24 20:19:50"
$ws.Range("E29").Value = "2026-02-24 20:19:53"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "89%"
$ws.Range("E30").Value = "2026-02-24 20:19:55"
$ws.Range("E31").Value = "2026-02-24 20:19:58"
$ws.Range("E32").Value = "2026-02-24 20:20:01"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "67%"
$ws.Range("O32").Value = "7.5 °C"
$ws.Range("E33").Value = "2026-02-24 20:20:03"
$ws.Range("J33").Value = "1021.6 hPa"
$ws.Range("K33").Value = "15.5 MJ/m2"
$ws.Range("E34").Value = "2026-02-24 20:20:06"
$ws.Range("E35").Value = "2026-02-24 20:20:09"
$ws.Range("E36").Value = "2026-02-24 20:20:12"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "80%"
$ws.Range("J36").Value = "1019.9 hPa"
$ws.Range("O36").Value = "13.1 °C"
$ws.Range("E37").Value = "2026-02-24 20:20:14"
$ws.Range("O37").Value = "8.9 °C"
$ws.Range("E38").Value = "2026-02-24 20:20:17"
$ws.Range("E39").Value = "2026-02-24 20:20:20"
$ws.Range("N39").Value = "1.6 °C 19:51 TU"
$ws.Range("E40").Value = "2026-02-24 20:20:22"
$ws.Range("E41").Value = "2026-02-24 20:20:25"
$ws.Range("E42").Value = "2026-02-24 20:20:28"
$ws.Range("O42").Value = "11.3 °C"
$ws.Range("E43").Value = "2026-02-24 20:20:30"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "69%"
$ws.Range("E44").Value = "2026-02-24 20:20:33"
$ws.Range("E45").Value = "2026-02-24 20:20:35"
$ws.Range("H45").NumberFormat = "@"
$ws.Range("H45").Value = "43%"
$ws.Range("E46").Value = "2026-02-24 20:20:38"
$ws.Range("L46").Value = "15.5 km/h - 178º 19:45 TU"
